# #327 Ajout des profils d'acces
#
# 1. Metadata sheet: bump the "Date" value.
# 2. Elements sheet: swap the two "Mapping" columns (AK <-> AL), content,
#    header text and column widths, so that the "Specification metier"
#    mapping now comes before the "RIM Mapping" one.

$wb = $excel.ActiveWorkbook

# --- 1. Metadata!B8 (Date) -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- 2. Elements sheet: swap columns AK (37) and AL (38) -------------------
$els = $wb.Worksheets.Item("Elements")

# Swap the cell contents row by row (only rows that actually hold data in
# either of the two columns need touching - rows 2 and 4 are empty on both
# sides already, so leave them untouched).
foreach ($r in 1,3,5,6) {
    $akCell = $els.Cells.Item($r, 37)
    $alCell = $els.Cells.Item($r, 38)
    $akText = $akCell.Text
    $alText = $alCell.Text
    $akCell.Value = $alText
    $alCell.Value = $akText
}

# Swap the column widths that go along with the swapped content.
# (AK was the narrow "n/a" style column, AL the wide "Specification" one -
# the stored OOXML widths are swapped between the two columns.)
$els.Columns.Item(37).ColumnWidth = 75.91796875
$els.Columns.Item(38).ColumnWidth = 24.0859375
